{"js": "// Replace each two-digit multiplication expression with its new value.\n// Every source expression in this worksheet is unique, so a case-sensitive\n// exact-text search safely targets exactly one run each.\nconst replacements = [\n  [\"41\u00d754=\", \"57\u00d785=\"],\n  [\"60\u00d731=\", \"87\u00d764=\"],\n  [\"22\u00d779=\", \"59\u00d787=\"],\n  [\"29\u00d739=\", \"47\u00d736=\"],\n  [\"55\u00d736=\", \"14\u00d788=\"],\n  [\"30\u00d765=\", \"75\u00d799=\"],\n  [\"34\u00d786=\", \"52\u00d750=\"],\n  [\"68\u00d779=\", \"89\u00d786=\"],\n  [\"69\u00d713=\", \"49\u00d765=\"],\n  [\"59\u00d795=\", \"41\u00d785=\"],\n  [\"90\u00d796=\", \"65\u00d760=\"],\n  [\"83\u00d753=\", \"77\u00d796=\"],\n  [\"13\u00d752=\", \"63\u00d777=\"],\n  [\"39\u00d774=\", \"24\u00d751=\"],\n  [\"98\u00d794=\", \"99\u00d795=\"],\n  [\"59\u00d757=\", \"36\u00d750=\"],\n  [\"63\u00d792=\", \"40\u00d724=\"],\n  [\"78\u00d728=\", \"45\u00d725=\"],\n  [\"64\u00d739=\", \"73\u00d784=\"],\n  [\"50\u00d766=\", \"18\u00d734=\"],\n  [\"89\u00d768=\", \"95\u00d766=\"],\n  [\"94\u00d787=\", \"62\u00d756=\"],\n  [\"57\u00d793=\", \"70\u00d765=\"],\n  [\"79\u00d790=\", \"84\u00d732=\"],\n  [\"35\u00d757=\", \"97\u00d759=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its new value.\n# Every source expression in this worksheet is unique, so a case-sensitive\n# whole-document Find/Replace safely targets exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"41\u00d754=\", \"57\u00d785=\"),\n  @(\"60\u00d731=\", \"87\u00d764=\"),\n  @(\"22\u00d779=\", \"59\u00d787=\"),\n  @(\"29\u00d739=\", \"47\u00d736=\"),\n  @(\"55\u00d736=\", \"14\u00d788=\"),\n  @(\"30\u00d765=\", \"75\u00d799=\"),\n  @(\"34\u00d786=\", \"52\u00d750=\"),\n  @(\"68\u00d779=\", \"89\u00d786=\"),\n  @(\"69\u00d713=\", \"49\u00d765=\"),\n  @(\"59\u00d795=\", \"41\u00d785=\"),\n  @(\"90\u00d796=\", \"65\u00d760=\"),\n  @(\"83\u00d753=\", \"77\u00d796=\"),\n  @(\"13\u00d752=\", \"63\u00d777=\"),\n  @(\"39\u00d774=\", \"24\u00d751=\"),\n  @(\"98\u00d794=\", \"99\u00d795=\"),\n  @(\"59\u00d757=\", \"36\u00d750=\"),\n  @(\"63\u00d792=\", \"40\u00d724=\"),\n  @(\"78\u00d728=\", \"45\u00d725=\"),\n  @(\"64\u00d739=\", \"73\u00d784=\"),\n  @(\"50\u00d766=\", \"18\u00d734=\"),\n  @(\"89\u00d768=\", \"95\u00d766=\"),\n  @(\"94\u00d787=\", \"62\u00d756=\"),\n  @(\"57\u00d793=\", \"70\u00d765=\"),\n  @(\"79\u00d790=\", \"84\u00d732=\"),\n  @(\"35\u00d757=\", \"97\u00d759=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $null = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
